$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 6 formatting (date style on A, shared "Bag" string on N) down to row 7
$ws.Range("A6:N6").Copy($ws.Range("A7:N7"))

$ws.Range("A7").Value = 42607.886365740742
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 54
$ws.Range("D7").Value = 42
$ws.Range("E7").Value = 36
$ws.Range("F7").Value = 63
$ws.Range("G7").Value = 11243
$ws.Range("H7").Value = 24673
$ws.Range("I7").Value = 2799
$ws.Range("J7").Value = 363
$ws.Range("K7").Value = 280
$ws.Range("L7").Value = 11
$ws.Range("M7").Value = 19
$ws.Range("N7").Value = "Bag"
